$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "other": update two existing parameter values and append the new
# "class 2" (pv2 / bat2) parameter block in rows 11-18.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Range("B2").Value = 0.2
$wsOther.Range("B3").Value = 150

$wsOther.Range("A11").Value = "pv2_eff"
$wsOther.Range("B11").Value = 0.2

$wsOther.Range("A12").Value = "pv2_area"
$wsOther.Range("B12").Value = 1

$wsOther.Range("A13").Value = "bat2_c_rate_ch"
$wsOther.Range("B13").Value = 1

$wsOther.Range("A14").Value = "bat2_c_rate_dis"
$wsOther.Range("B14").Value = 1

$wsOther.Range("A15").Value = "bat2_ch_eff"
$wsOther.Range("B15").Value = 0.95

$wsOther.Range("A16").Value = "bat2_dis_eff"
$wsOther.Range("B16").Value = 0.95

$wsOther.Range("A17").Value = "bat2_starting_SOC"
$wsOther.Range("B17").Value = 0.7

$wsOther.Range("A18").Value = "bat2_E_max"
$wsOther.Range("B18").Value = 200

$wsOther.Range("F22").Select()

# ---------------------------------------------------------------------------
# Sheet "series": move selection, set page orientation to portrait.
# ---------------------------------------------------------------------------
$wsSeries = $wb.Worksheets.Item("series")
$wsSeries.PageSetup.Orientation = 1
$wsSeries.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet "elements": move selection, set page orientation to portrait, and
# stamp an (empty) underline-styled cell at E4 - mirrors the formatting
# already used for the "last touched" marker cells on the "other" sheet.
# ---------------------------------------------------------------------------
$wsElements = $wb.Worksheets.Item("elements")
$wsElements.PageSetup.Orientation = 1
$wsElements.Range("E4").Font.Underline = 2
$wsElements.Range("B2").Select()

# ---------------------------------------------------------------------------
# Restore the originally active sheet/tab ("conect") - none of its own
# content changed, only the workbook-level active-tab bookkeeping should be
# left exactly as it was before this edit.
# ---------------------------------------------------------------------------
$wsConect = $wb.Worksheets.Item("conect")
$wsConect.Activate()
